# Sprint-1 Project Backlog is done
#
# Adds the three remaining Sprint #1 backlog items (Inverse, Latex and
# RemoveRepetitiveWords transformers) together with their tasks to the
# "Sprint #1 Backlog" sheet - mirroring the layout already used for the
# first two backlog items on that sheet (rows 5-8 and 11-14) - and updates
# the recorded selections on the "Product Backlog" / "Sprint #1 Backlog"
# sheets to reflect where the user ended up scrolled to/selecting.

$wb = $excel.ActiveWorkbook
$wsProd = $wb.Worksheets.Item("Product Backlog")
$wsSprint = $wb.Worksheets.Item("Sprint #1 Backlog")

$xlTop = -4160
$xlCenter = -4108

function Set-BacklogItemCell($range, $text) {
    $range.Value = $text
    $range.VerticalAlignment = $xlTop
    $range.WrapText = $true
}

function Set-LabourCell($range, $text) {
    $range.Value = $text
    $range.HorizontalAlignment = $xlCenter
    $range.VerticalAlignment = $xlCenter
}

function Set-TaskCell($range, $text) {
    $range.Value = $text
    $range.VerticalAlignment = $xlTop
    $range.WrapText = $true
}

# --- Backlog item #3: Inverse transformer (rows 17-19) ---
Set-BacklogItemCell $wsSprint.Range("A17") "As a user, I can reverse the sequence of characters, retaining the case of letters in specific positions (inverse)"
Set-LabourCell       $wsSprint.Range("C17") "Medium"
Set-TaskCell         $wsSprint.Range("B18") "Task #1 Implement InverseTransformer class with logic to preserve case positions."
Set-TaskCell         $wsSprint.Range("B19") "Task #2 Write unit tests for edge cases (e.g., single character, empty string, mixed case)."

# --- Backlog item #6: Latex transformer (rows 22-23) ---
Set-BacklogItemCell $wsSprint.Range("A22") "As a user I can convert text to Latex supported format (special characters) -> (John Smith & Sons -> John Smith \& Sons)"
Set-LabourCell       $wsSprint.Range("C22") "Low"
Set-TaskCell         $wsSprint.Range("B23") "Task #1 Implement LatexTransformer class to handle & and $."

# --- Backlog item #8: RemoveRepetitiveWords transformer (rows 26-28) ---
Set-BacklogItemCell $wsSprint.Range("A26") "As a user, I can eliminate repetitive words in the immediate vicinity (Send me me a message -> Send me a message)"
Set-LabourCell       $wsSprint.Range("C26") "Medium"
Set-TaskCell         $wsSprint.Range("B27") "Task #1 Implement RemoveRepetitiveWordsTransformer class using regular expressions or a loop."
Set-TaskCell         $wsSprint.Range("B28") "Task #2 Test with multiple repetitions"

# --- Row heights, matching the ht used by the equivalent rows already on the sheet ---
$wsSprint.Rows.Item(17).RowHeight = 45
$wsSprint.Rows.Item(18).RowHeight = 30
$wsSprint.Rows.Item(19).RowHeight = 30
$wsSprint.Rows.Item(22).RowHeight = 45
$wsSprint.Rows.Item(23).RowHeight = 30
$wsSprint.Rows.Item(26).RowHeight = 45
$wsSprint.Rows.Item(27).RowHeight = 30
$wsSprint.Rows.Item(28).RowHeight = 15

# --- View selections ---
# "Product Backlog" sheet: scrolled down a bit, selection on A9.
$wsProd.Activate()
$wsProd.Range("A9").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1

# "Sprint #1 Backlog" is the tab that should stay the active one, so select
# its range last, putting the selection on B26 and scrolling so row 9 is on top.
$wsSprint.Activate()
$wsSprint.Range("B26").Select()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
